$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("T55")

# New values for B2:G11 (each row's data is shifted up by one from the
# previous state, with a brand-new final row of computed values appended).
$data = @(
    @(-0.005035358036557557, 0.3637990656943072,  0.1855351867275108,  0.4307379559865961, 0.4469673578304695, 14),
    @(-0.01898232632975465,  0.3165751929851393,  0.1363897707384557,  0.3693098573534908, 0.3838817902180699, 13),
    @(0.00780077805212256,   0.3794346132818944,  0.1847968220473418,  0.4298800088947401, 0.4489210943938488, 12),
    @(0.03532231998103826,   0.2564852750112934,  0.1004647918617584,  0.3169618145167623, 0.3303616777566779, 11),
    @(0.002000156080238219,  0.3449448151542904,  0.1826303857083173,  0.4273527649475516, 0.4504644332784307, 10),
    @(-0.05875859174690282,  0.3248866822092542,  0.1283271032654252,  0.3582277254281489, 0.3748117494910135, 9),
    @(-0.02790046359007027,  0.3120137538489823,  0.1287774372279015,  0.3588557331684997, 0.3824716629792676, 8),
    @(-0.03188629724616485,  0.2657389849834738,  0.1191175486163901,  0.3451341023665874, 0.3711930514458828, 7),
    @(-0.0881379738596985,   0.4006321887415465,  0.2481252247171507,  0.4981216966938408, 0.5370552175868598, 6),
    @(0.003999580097293598,  0.1842291207164398,  0.04734795102614385, 0.2175958433108129, 0.2432384488140937, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}

$wb.Save()
